{"js": "// Replace the 100 two-digit-multiplication problems in the table, in document\n// order (row-major, 20 rows x 5 columns), old -> new text.\nconst replacements = [\n  [\"15\u00d784=\", \"37\u00d755=\"],\n  [\"60\u00d776=\", \"85\u00d745=\"],\n  [\"44\u00d795=\", \"74\u00d786=\"],\n  [\"13\u00d713=\", \"27\u00d769=\"],\n  [\"11\u00d774=\", \"23\u00d737=\"],\n  [\"57\u00d769=\", \"13\u00d741=\"],\n  [\"96\u00d733=\", \"75\u00d761=\"],\n  [\"10\u00d754=\", \"34\u00d796=\"],\n  [\"98\u00d769=\", \"40\u00d725=\"],\n  [\"32\u00d769=\", \"63\u00d720=\"],\n  [\"43\u00d792=\", \"96\u00d799=\"],\n  [\"35\u00d714=\", \"56\u00d790=\"],\n  [\"79\u00d785=\", \"51\u00d742=\"],\n  [\"32\u00d729=\", \"28\u00d748=\"],\n  [\"35\u00d761=\", \"27\u00d725=\"],\n  [\"39\u00d734=\", \"10\u00d748=\"],\n  [\"86\u00d712=\", \"69\u00d741=\"],\n  [\"74\u00d755=\", \"61\u00d772=\"],\n  [\"55\u00d768=\", \"20\u00d772=\"],\n  [\"12\u00d788=\", \"10\u00d738=\"],\n  [\"79\u00d773=\", \"38\u00d716=\"],\n  [\"16\u00d798=\", \"42\u00d723=\"],\n  [\"22\u00d781=\", \"61\u00d747=\"],\n  [\"60\u00d788=\", \"22\u00d723=\"],\n  [\"94\u00d764=\", \"66\u00d785=\"],\n  [\"90\u00d763=\", \"68\u00d711=\"],\n  [\"90\u00d769=\", \"82\u00d774=\"],\n  [\"55\u00d769=\", \"69\u00d737=\"],\n  [\"84\u00d781=\", \"40\u00d713=\"],\n  [\"10\u00d763=\", \"56\u00d721=\"],\n  [\"12\u00d724=\", \"24\u00d763=\"],\n  [\"20\u00d7100=\", \"33\u00d750=\"],\n  [\"28\u00d758=\", \"25\u00d736=\"],\n  [\"65\u00d710=\", \"61\u00d752=\"],\n  [\"30\u00d721=\", \"29\u00d755=\"],\n  [\"14\u00d744=\", \"13\u00d762=\"],\n  [\"16\u00d762=\", \"95\u00d775=\"],\n  [\"86\u00d788=\", \"52\u00d797=\"],\n  [\"19\u00d797=\", \"13\u00d776=\"],\n  [\"23\u00d714=\", \"87\u00d720=\"],\n  [\"71\u00d743=\", \"90\u00d737=\"],\n  [\"53\u00d781=\", \"28\u00d790=\"],\n  [\"19\u00d791=\", \"48\u00d716=\"],\n  [\"67\u00d750=\", \"94\u00d776=\"],\n  [\"32\u00d751=\", \"34\u00d731=\"],\n  [\"89\u00d729=\", \"25\u00d787=\"],\n  [\"12\u00d722=\", \"61\u00d731=\"],\n  [\"36\u00d732=\", \"45\u00d769=\"],\n  [\"57\u00d781=\", \"87\u00d780=\"],\n  [\"19\u00d798=\", \"36\u00d753=\"],\n  [\"39\u00d777=\", \"52\u00d721=\"],\n  [\"72\u00d730=\", \"42\u00d757=\"],\n  [\"45\u00d795=\", \"81\u00d798=\"],\n  [\"11\u00d777=\", \"55\u00d790=\"],\n  [\"90\u00d773=\", \"30\u00d785=\"],\n  [\"81\u00d747=\", \"33\u00d757=\"],\n  [\"36\u00d711=\", \"77\u00d788=\"],\n  [\"32\u00d791=\", \"48\u00d750=\"],\n  [\"20\u00d792=\", \"95\u00d781=\"],\n  [\"80\u00d785=\", \"86\u00d730=\"],\n  [\"62\u00d744=\", \"23\u00d794=\"],\n  [\"10\u00d740=\", \"48\u00d777=\"],\n  [\"76\u00d739=\", \"11\u00d720=\"],\n  [\"32\u00d752=\", \"77\u00d767=\"],\n  [\"13\u00d752=\", \"75\u00d731=\"],\n  [\"50\u00d730=\", \"65\u00d774=\"],\n  [\"58\u00d724=\", \"58\u00d758=\"],\n  [\"56\u00d797=\", \"28\u00d791=\"],\n  [\"68\u00d767=\", \"87\u00d751=\"],\n  [\"98\u00d793=\", \"41\u00d788=\"],\n  [\"37\u00d720=\", \"28\u00d788=\"],\n  [\"86\u00d798=\", \"13\u00d764=\"],\n  [\"17\u00d718=\", \"46\u00d751=\"],\n  [\"81\u00d712=\", \"92\u00d754=\"],\n  [\"99\u00d786=\", \"17\u00d752=\"],\n  [\"83\u00d798=\", \"92\u00d790=\"],\n  [\"15\u00d772=\", \"43\u00d762=\"],\n  [\"32\u00d786=\", \"89\u00d739=\"],\n  [\"59\u00d796=\", \"100\u00d768=\"],\n  [\"63\u00d788=\", \"95\u00d751=\"],\n  [\"77\u00d774=\", \"44\u00d748=\"],\n  [\"49\u00d794=\", \"18\u00d720=\"],\n  [\"16\u00d794=\", \"48\u00d757=\"],\n  [\"68\u00d730=\", \"90\u00d765=\"],\n  [\"46\u00d784=\", \"34\u00d725=\"],\n  [\"49\u00d728=\", \"38\u00d718=\"],\n  [\"26\u00d725=\", \"46\u00d774=\"],\n  [\"76\u00d793=\", \"90\u00d732=\"],\n  [\"68\u00d763=\", \"51\u00d782=\"],\n  [\"88\u00d734=\", \"85\u00d728=\"],\n  [\"24\u00d798=\", \"23\u00d787=\"],\n  [\"59\u00d727=\", \"35\u00d774=\"],\n  [\"45\u00d764=\", \"18\u00d796=\"],\n  [\"18\u00d724=\", \"86\u00d764=\"],\n  [\"60\u00d780=\", \"27\u00d734=\"],\n  [\"48\u00d796=\", \"16\u00d723=\"],\n  [\"59\u00d797=\", \"73\u00d752=\"],\n  [\"44\u00d722=\", \"14\u00d723=\"],\n  [\"86\u00d760=\", \"45\u00d769=\"],\n  [\"28\u00d755=\", \"77\u00d781=\"]\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\nconst columns = 5;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount && idx < replacements.length; r++) {\n  for (let c = 0; c < columns && idx < replacements.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const [oldText, newText] = replacements[idx];\n    const para = cell.body.paragraphs.items[0];\n    para.load(\"text\");\n    await context.sync();\n\n    if (para.text !== oldText) {\n      throw new Error(\n        \"Cell (\" + r + \",\" + c + \") text mismatch: expected '\" + oldText +\n        \"' but found '\" + para.text + \"'\"\n      );\n    }\n\n    para.insertText(newText, Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 two-digit-multiplication problems in the table, in\n# document order (row-major, 20 rows x 5 columns), old -> new text.\n$replacements = @(\n    @(\"15\u00d784=\", \"37\u00d755=\"),\n    @(\"60\u00d776=\", \"85\u00d745=\"),\n    @(\"44\u00d795=\", \"74\u00d786=\"),\n    @(\"13\u00d713=\", \"27\u00d769=\"),\n    @(\"11\u00d774=\", \"23\u00d737=\"),\n    @(\"57\u00d769=\", \"13\u00d741=\"),\n    @(\"96\u00d733=\", \"75\u00d761=\"),\n    @(\"10\u00d754=\", \"34\u00d796=\"),\n    @(\"98\u00d769=\", \"40\u00d725=\"),\n    @(\"32\u00d769=\", \"63\u00d720=\"),\n    @(\"43\u00d792=\", \"96\u00d799=\"),\n    @(\"35\u00d714=\", \"56\u00d790=\"),\n    @(\"79\u00d785=\", \"51\u00d742=\"),\n    @(\"32\u00d729=\", \"28\u00d748=\"),\n    @(\"35\u00d761=\", \"27\u00d725=\"),\n    @(\"39\u00d734=\", \"10\u00d748=\"),\n    @(\"86\u00d712=\", \"69\u00d741=\"),\n    @(\"74\u00d755=\", \"61\u00d772=\"),\n    @(\"55\u00d768=\", \"20\u00d772=\"),\n    @(\"12\u00d788=\", \"10\u00d738=\"),\n    @(\"79\u00d773=\", \"38\u00d716=\"),\n    @(\"16\u00d798=\", \"42\u00d723=\"),\n    @(\"22\u00d781=\", \"61\u00d747=\"),\n    @(\"60\u00d788=\", \"22\u00d723=\"),\n    @(\"94\u00d764=\", \"66\u00d785=\"),\n    @(\"90\u00d763=\", \"68\u00d711=\"),\n    @(\"90\u00d769=\", \"82\u00d774=\"),\n    @(\"55\u00d769=\", \"69\u00d737=\"),\n    @(\"84\u00d781=\", \"40\u00d713=\"),\n    @(\"10\u00d763=\", \"56\u00d721=\"),\n    @(\"12\u00d724=\", \"24\u00d763=\"),\n    @(\"20\u00d7100=\", \"33\u00d750=\"),\n    @(\"28\u00d758=\", \"25\u00d736=\"),\n    @(\"65\u00d710=\", \"61\u00d752=\"),\n    @(\"30\u00d721=\", \"29\u00d755=\"),\n    @(\"14\u00d744=\", \"13\u00d762=\"),\n    @(\"16\u00d762=\", \"95\u00d775=\"),\n    @(\"86\u00d788=\", \"52\u00d797=\"),\n    @(\"19\u00d797=\", \"13\u00d776=\"),\n    @(\"23\u00d714=\", \"87\u00d720=\"),\n    @(\"71\u00d743=\", \"90\u00d737=\"),\n    @(\"53\u00d781=\", \"28\u00d790=\"),\n    @(\"19\u00d791=\", \"48\u00d716=\"),\n    @(\"67\u00d750=\", \"94\u00d776=\"),\n    @(\"32\u00d751=\", \"34\u00d731=\"),\n    @(\"89\u00d729=\", \"25\u00d787=\"),\n    @(\"12\u00d722=\", \"61\u00d731=\"),\n    @(\"36\u00d732=\", \"45\u00d769=\"),\n    @(\"57\u00d781=\", \"87\u00d780=\"),\n    @(\"19\u00d798=\", \"36\u00d753=\"),\n    @(\"39\u00d777=\", \"52\u00d721=\"),\n    @(\"72\u00d730=\", \"42\u00d757=\"),\n    @(\"45\u00d795=\", \"81\u00d798=\"),\n    @(\"11\u00d777=\", \"55\u00d790=\"),\n    @(\"90\u00d773=\", \"30\u00d785=\"),\n    @(\"81\u00d747=\", \"33\u00d757=\"),\n    @(\"36\u00d711=\", \"77\u00d788=\"),\n    @(\"32\u00d791=\", \"48\u00d750=\"),\n    @(\"20\u00d792=\", \"95\u00d781=\"),\n    @(\"80\u00d785=\", \"86\u00d730=\"),\n    @(\"62\u00d744=\", \"23\u00d794=\"),\n    @(\"10\u00d740=\", \"48\u00d777=\"),\n    @(\"76\u00d739=\", \"11\u00d720=\"),\n    @(\"32\u00d752=\", \"77\u00d767=\"),\n    @(\"13\u00d752=\", \"75\u00d731=\"),\n    @(\"50\u00d730=\", \"65\u00d774=\"),\n    @(\"58\u00d724=\", \"58\u00d758=\"),\n    @(\"56\u00d797=\", \"28\u00d791=\"),\n    @(\"68\u00d767=\", \"87\u00d751=\"),\n    @(\"98\u00d793=\", \"41\u00d788=\"),\n    @(\"37\u00d720=\", \"28\u00d788=\"),\n    @(\"86\u00d798=\", \"13\u00d764=\"),\n    @(\"17\u00d718=\", \"46\u00d751=\"),\n    @(\"81\u00d712=\", \"92\u00d754=\"),\n    @(\"99\u00d786=\", \"17\u00d752=\"),\n    @(\"83\u00d798=\", \"92\u00d790=\"),\n    @(\"15\u00d772=\", \"43\u00d762=\"),\n    @(\"32\u00d786=\", \"89\u00d739=\"),\n    @(\"59\u00d796=\", \"100\u00d768=\"),\n    @(\"63\u00d788=\", \"95\u00d751=\"),\n    @(\"77\u00d774=\", \"44\u00d748=\"),\n    @(\"49\u00d794=\", \"18\u00d720=\"),\n    @(\"16\u00d794=\", \"48\u00d757=\"),\n    @(\"68\u00d730=\", \"90\u00d765=\"),\n    @(\"46\u00d784=\", \"34\u00d725=\"),\n    @(\"49\u00d728=\", \"38\u00d718=\"),\n    @(\"26\u00d725=\", \"46\u00d774=\"),\n    @(\"76\u00d793=\", \"90\u00d732=\"),\n    @(\"68\u00d763=\", \"51\u00d782=\"),\n    @(\"88\u00d734=\", \"85\u00d728=\"),\n    @(\"24\u00d798=\", \"23\u00d787=\"),\n    @(\"59\u00d727=\", \"35\u00d774=\"),\n    @(\"45\u00d764=\", \"18\u00d796=\"),\n    @(\"18\u00d724=\", \"86\u00d764=\"),\n    @(\"60\u00d780=\", \"27\u00d734=\"),\n    @(\"48\u00d796=\", \"16\u00d723=\"),\n    @(\"59\u00d797=\", \"73\u00d752=\"),\n    @(\"44\u00d722=\", \"14\u00d723=\"),\n    @(\"86\u00d760=\", \"45\u00d769=\"),\n    @(\"28\u00d755=\", \"77\u00d781=\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$columns = 5\n$idx = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $columns; $c++) {\n        if ($idx -ge $replacements.Count) { break }\n\n        $pair = $replacements[$idx]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n\n        $cell = $tbl.Cell($r, $c)\n        $cellRange = $cell.Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        $current = $cellRange.Text\n\n        if ($current -ne $oldText) {\n            throw (\"Cell (\" + $r + \",\" + $c + \") text mismatch: expected '\" + $oldText + \"' but found '\" + $current + \"'\")\n        }\n\n        $cellRange.Text = $newText\n        $idx++\n    }\n}\n"}
